$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "67.066.69"
$ws.Range("E2").Value = "  +0.02%  "
$ws.Range("D3").Value = "3.093.47"
$ws.Range("E3").Value = "  +0.25%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Formula = "'572.14"
$ws.Range("E5").Value = "  -1.12%  "
$ws.Range("D6").Formula = "'177.08"
$ws.Range("E6").Value = "  +4.89%  "
$ws.Range("E7").Value = "  +0.11%  "
$ws.Range("D8").Value = "3.090.17"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").Formula = "'0.513"
$ws.Range("E9").Value = "  -0.43%  "
$ws.Range("D10").Formula = "'6.41"
$ws.Range("E10").Value = "  -0.28%  "
$ws.Range("D11").Formula = "'0.152"
$ws.Range("E11").Value = "  +0.98%  "
$ws.Range("D12").Formula = "'0.468"
$ws.Range("E12").Value = "  -1.00%  "
$ws.Range("E13").Value = "  -0.28%  "
$ws.Range("D14").Formula = "'35.94"
$ws.Range("E14").Value = "  -0.69%  "
$ws.Range("E15").Value = "  +0.80%  "
$ws.Range("D16").Value = "3.609.93"
$ws.Range("E16").Value = "  +0.39%  "
$ws.Range("D17").Value = "67.045.92"
$ws.Range("E17").Value = "  +0.05%  "
$ws.Range("D18").Formula = "'7.02"
$ws.Range("E18").Value = "  +0.13%  "
$ws.Range("D19").Value = "3.088.89"
$ws.Range("E19").Value = "  +0.05%  "
$ws.Range("D20").Formula = "'16.52"
$ws.Range("E20").Value = "  -0.62%  "
$ws.Range("D21").Formula = "'486.51"
$ws.Range("E21").Value = "  -0.87%  "
$ws.Range("D22").Formula = "'7.70"
$ws.Range("E22").Value = "  -0.88%  "
$ws.Range("D23").Formula = "'0.685"
$ws.Range("E23").Value = "  -0.39%  "
$ws.Range("D24").Formula = "'83.40"
$ws.Range("E24").Value = "  +0.73%  "
$ws.Range("D25").Formula = "'12.78"
$ws.Range("E25").Value = "  -0.90%  "
$ws.Range("E26").Value = "  +0.65%  "
$ws.Range("D27").Formula = "'10.21"
$ws.Range("E27").Value = "  -0.48%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("D29").Formula = "'7.88"
$ws.Range("E29").Value = "  +0.77%  "
$ws.Range("D30").Formula = "'2.30"
$ws.Range("E30").Value = "  -0.91%  "
$ws.Range("D31").Formula = "'2.58"
$ws.Range("E31").Value = "  -2.05%  "
$ws.Range("D32").Formula = "'28.01"
$ws.Range("E32").Value = "  +0.40%  "
$ws.Range("D33").Formula = "'0.112"
$ws.Range("E33").Value = "  -0.08%  "
$ws.Range("D34").Value = "0.0₃0941"
$ws.Range("E34").Value = "  +2.99%  "
$ws.Range("D35").Formula = "'0.998"
$ws.Range("E35").Value = "  -0.30%  "
$ws.Range("D36").Formula = "'47.24"
$ws.Range("E36").Value = "  +1.76%  "
$ws.Range("D37").Formula = "'0.947"
$ws.Range("E37").Value = "  -0.78%  "
$ws.Range("E38").Value = "  -2.47%  "
$ws.Range("D39").Formula = "'0.311"
$ws.Range("E39").Value = "  +2.75%  "
$ws.Range("B40").Value = "Stacks"
$ws.Range("C40").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D40").Formula = "'2.01"
$ws.Range("E40").Value = "  +0.85%  "
$ws.Range("B41").Value = "OKB"
$ws.Range("C41").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D41").Formula = "'48.92"
$ws.Range("E41").Value = "  -0.74%  "
$ws.Range("E42").Value = "  -0.38%  "
$ws.Range("D43").Formula = "'2.72"
$ws.Range("E43").Value = "  +9.33%  "
$ws.Range("E44").Value = "  -1.49%  "
$ws.Range("D45").Value = "2.807.86"
$ws.Range("E45").Value = "  +1.00%  "
$ws.Range("D46").Formula = "'369.54"
$ws.Range("E46").Value = "  -0.75%  "
$ws.Range("D47").Formula = "'0.0344"
$ws.Range("E47").Value = "  -0.47%  "
$ws.Range("D48").Formula = "'134.57"
$ws.Range("E48").Value = "  -0.72%  "
$ws.Range("E49").Value = "  +0.04%  "
$ws.Range("D50").Formula = "'25.63"
$ws.Range("E50").Value = "  +4.45%  "
$ws.Range("E51").Value = "  +6.85%  "
